# Update the team-specific transition matrix on Sheet1 with the latest
# simulated-game results (more games simulated; several previously-empty
# transition probabilities are now populated).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (state 0)
$ws.Range("B2").Value = 0.05263157894736842
$ws.Range("C2").Value = 0.8421052631578947
$ws.Range("P2").Value = 0.1052631578947368

# Row 3 (state 1)
$ws.Range("P3").Value = 0.5
$ws.Range("S3").Value = 0.5

# Row 4 (state 2)
$ws.Range("P4").Value = 0.5
$ws.Range("S4").Value = 0.5

# Row 6 (state 4)
$ws.Range("J6").Value = 0.1
$ws.Range("Q6").Value = 0.2
$ws.Range("S6").Value = 0.7

# Row 7 (state 5)
$ws.Range("F7").Value = 0.1
$ws.Range("J7").Value = 0.3
$ws.Range("Q7").Value = 0.1
$ws.Range("S7").Value = 0.5

# Row 8 (state 6)
$ws.Range("B8").Value = 0.08108108108108109
$ws.Range("F8").Value = 0.02702702702702703
$ws.Range("J8").Value = 0.2432432432432433
$ws.Range("O8").Value = 0.05405405405405406
$ws.Range("Q8").Value = 0.1081081081081081
$ws.Range("R8").Value = 0.1351351351351351
$ws.Range("S8").Value = 0.3513513513513514

# Row 9 (state 7)
$ws.Range("B9").Value = 0.1
$ws.Range("F9").Value = 0.05
$ws.Range("J9").Value = 0.15
$ws.Range("Q9").Value = 0.15
$ws.Range("R9").Value = 0.1
$ws.Range("S9").Value = 0.45

# Row 10 (state 8)
$ws.Range("B10").Value = 0.1558441558441558
$ws.Range("D10").Value = 0.02597402597402598
$ws.Range("F10").Value = 0.05194805194805195
$ws.Range("J10").Value = 0.1038961038961039
$ws.Range("O10").Value = 0.02597402597402598
$ws.Range("Q10").Value = 0.1818181818181818
$ws.Range("R10").Value = 0.06493506493506493
$ws.Range("S10").Value = 0.3896103896103896

# Row 11 (state 9)
$ws.Range("G11").Value = 0.1052631578947368
$ws.Range("J11").Value = 0.1052631578947368
$ws.Range("K11").Value = 0.2631578947368421
$ws.Range("L11").Value = 0.4736842105263158
$ws.Range("S11").Value = 0.05263157894736842

# Row 12 (state 10)
$ws.Range("G12").Value = 0.6666666666666666
$ws.Range("J12").Value = 0.1111111111111111
$ws.Range("S12").Value = 0.2222222222222222

# Row 13 (state 11)
$ws.Range("G13").Value = 0.5
$ws.Range("J13").Value = 0.5

# Row 15 (state 13)
$ws.Range("H15").Value = 0.5
$ws.Range("J15").Value = 0.2142857142857143
$ws.Range("S15").Value = 0.2857142857142857

# Row 16 (state 14)
$ws.Range("H16").Value = 0.2727272727272727
$ws.Range("J16").Value = 0.5454545454545454
$ws.Range("S16").Value = 0.1818181818181818

# Row 17 (state 15)
$ws.Range("H17").Value = 0.04166666666666666
$ws.Range("I17").Value = 0.125
$ws.Range("J17").Value = 0.2916666666666667
$ws.Range("K17").Value = 0.1666666666666667
$ws.Range("M17").Value = 0.08333333333333333
$ws.Range("O17").Value = 0.08333333333333333
$ws.Range("S17").Value = 0.2083333333333333

# Row 18 (state 16)
$ws.Range("H18").Value = 0.08333333333333333
$ws.Range("I18").Value = 0.25
$ws.Range("J18").Value = 0.5
$ws.Range("O18").Value = 0.08333333333333333
$ws.Range("S18").Value = 0.08333333333333333

# Row 19 (state 17)
$ws.Range("F19").Value = 0.009523809523809525
$ws.Range("H19").Value = 0.2476190476190476
$ws.Range("I19").Value = 0.1333333333333333
$ws.Range("J19").Value = 0.2761904761904762
$ws.Range("K19").Value = 0.09523809523809523
$ws.Range("M19").Value = 0.02857142857142857
$ws.Range("O19").Value = 0.04761904761904762
$ws.Range("S19").Value = 0.1619047619047619
